# Auto-generated edit script applying the Ragnarok_Profits.xlsx diff
# Updates currentAveragePrice / LevePriceNQ / LeveProfitNQ / LeveProfitHQ style
# columns (H, I, J, K, L, M, N) on specific rows across 8 sheets, matching a
# scheduled-runner refresh of market-price snapshots.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 12151.869
$ws.Range("I69").Value = 9049.6
$ws.Range("K69").Value = 27148.8
$ws.Range("M69").Value = -26274.8
$ws.Range("H72").Value = 12151.869
$ws.Range("I72").Value = 9049.6
$ws.Range("K72").Value = 81446.40000000001
$ws.Range("M72").Value = -77078.40000000001
$ws.Range("H82").Value = 8560
$ws.Range("I82").Value = 6138.8335
$ws.Range("K82").Value = 18416.5005
$ws.Range("M82").Value = -18010.5005
$ws.Range("H85").Value = 8560
$ws.Range("I85").Value = 6138.8335
$ws.Range("K85").Value = 18416.5005
$ws.Range("M85").Value = -17012.5005
$ws.Range("H100").Value = 5658.6875
$ws.Range("I100").Value = 4191.2
$ws.Range("J100").Value = 8104.5
$ws.Range("K100").Value = 4191.2
$ws.Range("L100").Value = 8104.5
$ws.Range("M100").Value = -3650.2
$ws.Range("N100").Value = -9186.5
$ws.Range("H103").Value = 38463280
$ws.Range("I103").Value = 749.5
$ws.Range("J103").Value = 55557740
$ws.Range("K103").Value = 2248.5
$ws.Range("L103").Value = 166673220
$ws.Range("M103").Value = -1662.5
$ws.Range("N103").Value = -166674392
$ws.Range("H112").Value = 4137.6665
$ws.Range("J112").Value = 4322.2354
$ws.Range("L112").Value = 12966.7062
$ws.Range("N112").Value = -15182.7062
$ws.Range("H132").Value = 2309.6553
$ws.Range("I132").Value = 1971.24
$ws.Range("J132").Value = 4424.75
$ws.Range("K132").Value = 5913.72
$ws.Range("L132").Value = 13274.25
$ws.Range("M132").Value = -3383.72
$ws.Range("N132").Value = -18334.25
$ws.Range("H137").Value = 2000.2941
$ws.Range("I137").Value = 1664.9131
$ws.Range("J137").Value = 2701.5454
$ws.Range("K137").Value = 4994.7393
$ws.Range("L137").Value = 8104.6362
$ws.Range("M137").Value = -2444.7393
$ws.Range("N137").Value = -13204.6362
$ws.Range("H138").Value = 4065.2153
$ws.Range("J138").Value = 5067.3096
$ws.Range("L138").Value = 15201.9288
$ws.Range("N138").Value = -25481.9288
$ws.Range("H141").Value = 5982.933
$ws.Range("I141").Value = 6187.3335
$ws.Range("K141").Value = 18562.0005
$ws.Range("M141").Value = -13382.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 795.1667
$ws.Range("J2").Value = 880
$ws.Range("L2").Value = 880
$ws.Range("N2").Value = -1106
$ws.Range("H32").Value = 10014.293
$ws.Range("I32").Value = 9151.49
$ws.Range("K32").Value = 9151.49
$ws.Range("M32").Value = -8864.49
$ws.Range("H116").Value = 795.1667
$ws.Range("J116").Value = 880
$ws.Range("L116").Value = 880
$ws.Range("N116").Value = -5468
$ws.Range("H132").Value = 2129727.5
$ws.Range("I132").Value = 1962.409
$ws.Range("K132").Value = 5887.227000000001
$ws.Range("M132").Value = -3357.227000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 795.1667
$ws.Range("J3").Value = 880
$ws.Range("L3").Value = 880
$ws.Range("N3").Value = -1108
$ws.Range("H26").Value = 25964.5
$ws.Range("I26").Value = 10947
$ws.Range("K26").Value = 10947
$ws.Range("M26").Value = -10655
$ws.Range("H105").Value = 481021.62
$ws.Range("I105").Value = 759664.3
$ws.Range("K105").Value = 759664.3
$ws.Range("M105").Value = -757917.3
$ws.Range("H134").Value = 4003679.5
$ws.Range("I134").Value = 3698.6191
$ws.Range("J134").Value = 25003578
$ws.Range("K134").Value = 11095.8573
$ws.Range("L134").Value = 75010734
$ws.Range("M134").Value = -8560.8573
$ws.Range("N134").Value = -75015804

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 892.93335
$ws.Range("I16").Value = 787.25
$ws.Range("K16").Value = 787.25
$ws.Range("M16").Value = -500.25
$ws.Range("H31").Value = 29415278
$ws.Range("I31").Value = 50003012
$ws.Range("K31").Value = 50003012
$ws.Range("M31").Value = -50002717
$ws.Range("H34").Value = 29415278
$ws.Range("I34").Value = 50003012
$ws.Range("K34").Value = 50003012
$ws.Range("M34").Value = -50002810
$ws.Range("H58").Value = 2318.2222
$ws.Range("I58").Value = 1759.579
$ws.Range("K58").Value = 1759.579
$ws.Range("M58").Value = -1556.579
$ws.Range("H59").Value = 132500
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H69").Value = 52499.5
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 5000
$ws.Range("M69").Value = -4251
$ws.Range("H72").Value = 52499.5
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 15000
$ws.Range("M72").Value = -11256
$ws.Range("H113").Value = 892.93335
$ws.Range("I113").Value = 787.25
$ws.Range("K113").Value = 787.25
$ws.Range("M113").Value = 1382.75
$ws.Range("H132").Value = 1817.3478
$ws.Range("I132").Value = 1905.05
$ws.Range("J132").Value = 1232.6666
$ws.Range("K132").Value = 5715.15
$ws.Range("L132").Value = 3697.9998
$ws.Range("M132").Value = -3185.15
$ws.Range("N132").Value = -8757.9998
$ws.Range("H134").Value = 2034.919
$ws.Range("I134").Value = 1824.4286
$ws.Range("K134").Value = 5473.2858
$ws.Range("M134").Value = -2938.2858
$ws.Range("H136").Value = 2318.2222
$ws.Range("I136").Value = 1759.579
$ws.Range("K136").Value = 5278.737
$ws.Range("M136").Value = -2728.737

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4490.0454
$ws.Range("I131").Value = 2104.9092
$ws.Range("K131").Value = 6314.7276
$ws.Range("M131").Value = -1274.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 5692.1665
$ws.Range("J13").Value = 8199.75
$ws.Range("L13").Value = 8199.75
$ws.Range("N13").Value = -8477.75
$ws.Range("H62").Value = 42749.75
$ws.Range("J62").Value = 49333
$ws.Range("L62").Value = 49333
$ws.Range("N62").Value = -50705
$ws.Range("H65").Value = 42749.75
$ws.Range("J65").Value = 49333
$ws.Range("L65").Value = 147999
$ws.Range("N65").Value = -154863
$ws.Range("H128").Value = 99991.63
$ws.Range("J128").Value = 99991.63
$ws.Range("L128").Value = 99991.63
$ws.Range("N128").Value = -109951.63
$ws.Range("H132").Value = 1821649.1
$ws.Range("I132").Value = 3454.681
$ws.Range("K132").Value = 10364.043
$ws.Range("M132").Value = -7834.043

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 946.95
$ws.Range("I46").Value = 498.35715
$ws.Range("K46").Value = 498.35715
$ws.Range("M46").Value = -310.35715
$ws.Range("H74").Value = 94999.5
$ws.Range("I74").Value = 90000
$ws.Range("K74").Value = 90000
$ws.Range("M74").Value = -89002
$ws.Range("H77").Value = 94999.5
$ws.Range("I77").Value = 90000
$ws.Range("K77").Value = 270000
$ws.Range("M77").Value = -265008
$ws.Range("H93").Value = 3272112
$ws.Range("I93").Value = 3020.1667
$ws.Range("K93").Value = 3020.1667
$ws.Range("M93").Value = -1772.1667
$ws.Range("H122").Value = 3856.4893
$ws.Range("I122").Value = 3404.561
$ws.Range("K122").Value = 10213.683
$ws.Range("M122").Value = -7763.683000000001
$ws.Range("H130").Value = 98399
$ws.Range("J130").Value = 98399
$ws.Range("L130").Value = 98399
$ws.Range("N130").Value = -108439
$ws.Range("H136").Value = 4109.7407
$ws.Range("I136").Value = 2468.7058
$ws.Range("K136").Value = 7406.117400000001
$ws.Range("M136").Value = -4856.117400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2675.3157
$ws.Range("I122").Value = 1735.5333
$ws.Range("J122").Value = 6199.5
$ws.Range("K122").Value = 5206.5999
$ws.Range("L122").Value = 18598.5
$ws.Range("M122").Value = -2756.5999
$ws.Range("N122").Value = -23498.5
$ws.Range("H132").Value = 559454.3
$ws.Range("I132").Value = 4183.963
$ws.Range("K132").Value = 12551.889
$ws.Range("M132").Value = -10021.889
$ws.Range("H136").Value = 324617.94
$ws.Range("I136").Value = 2248.6191
$ws.Range("K136").Value = 6745.8573
$ws.Range("M136").Value = -4195.8573
$ws.Range("H137").Value = 34999
$ws.Range("J137").Value = 34999
$ws.Range("L137").Value = 34999
$ws.Range("N137").Value = -45199
